# Fruta / hortaliza, semanal
# A new weekly price record is inserted at the top of the data table
# (row 20), pushing the existing rows 20-28 down to 21-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20 (shifts rows 20:28 -> 21:29,
# carrying the D-column date style down with them).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly record.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44785
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112013
$ws.Range("G20").Value = "Alcachofa"
$ws.Range("H20").Value = "Argentina(o)"
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("N20").Value = "$/caja 50 unidades"
$ws.Range("O20").Value = "Región de Coquimbo"
$ws.Range("P20").Value = 310
$ws.Range("Q20").Value = 50
$ws.Range("R20").Value = "Hortaliza"
